# Update market_health_data.xlsx per 2025-10-30 10:34 data refresh

$wb = $excel.ActiveWorkbook

# --- Update "Metadata" sheet timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "30 Oct 2025, 10:34 AM"

# --- Update "distance from Dma50" sheet values (column C, rows 2-30) ---
$wsDma = $wb.Worksheets.Item("distance from Dma50")

$wsDma.Range("C2").Value = 9.9756
$wsDma.Range("C3").Value = 7.6922
$wsDma.Range("C4").Value = 6.2259
$wsDma.Range("C5").Value = 5.2519
$wsDma.Range("C6").Value = 5.0603
$wsDma.Range("C7").Value = 4.8612
$wsDma.Range("C8").Value = 4.6051
$wsDma.Range("C9").Value = 4.4789
$wsDma.Range("C10").Value = 3.6554
$wsDma.Range("C11").Value = 3.6367
$wsDma.Range("C12").Value = 3.4421
$wsDma.Range("C13").Value = 3.2058
$wsDma.Range("C14").Value = 3.2023
$wsDma.Range("C15").Value = 3.0995
$wsDma.Range("C16").Value = 3.0702
$wsDma.Range("C17").Value = 2.8566
$wsDma.Range("C18").Value = 2.5912
$wsDma.Range("C19").Value = 2.4892
$wsDma.Range("C20").Value = 2.3655
$wsDma.Range("C21").Value = 2.3187
$wsDma.Range("C22").Value = 1.3982
$wsDma.Range("C23").Value = 1.3686
$wsDma.Range("C24").Value = 1.2154
$wsDma.Range("C25").Value = 1.075
$wsDma.Range("C26").Value = 1.0004
$wsDma.Range("C27").Value = 0.9092
$wsDma.Range("C28").Value = 0.6111
$wsDma.Range("C29").Value = -0.0585
$wsDma.Range("C30").Value = -2.1298
